$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.629.79'
$ws.Range("E2").Value = '  -1.17%  '

# Row 3
$ws.Range("D3").Value = '1.851.59'
$ws.Range("E3").Value = '  -2.13%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  -0.87%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.41'
$ws.Range("E5").Value = '  -0.55%  '

# Row 6
$ws.Range("E6").Value = '  -0.83%  '

# Row 7
$ws.Range("E7").Value = '  -1.68%  '

# Row 8
$ws.Range("E8").Value = '  -0.70%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.03'
$ws.Range("E9").Value = '  -2.85%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07934'
$ws.Range("E10").Value = '  -1.06%  '

# Row 11
$ws.Range("E11").Value = '  -1.65%  '

# Row 12
$ws.Range("E12").Value = '  -1.04%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.959'
$ws.Range("E13").Value = '  +0.03%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.844.42'
$ws.Range("E14").Value = '  -2.41%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.170'
$ws.Range("E15").Value = '  +0.99%  '

# Row 16
$ws.Range("E16").Value = '  -0.92%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.36'
$ws.Range("E17").Value = '  +1.15%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06705'
$ws.Range("E18").Value = '  -1.08%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001035'
$ws.Range("E19").Value = '  -1.39%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.23'
$ws.Range("E20").Value = '  +0.68%  '

# Row 21
$ws.Range("E21").Value = '  -0.80%  '

# Row 22
$ws.Range("D22").Value = '27.615.14'
$ws.Range("E22").Value = '  -1.30%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.407'
$ws.Range("E23").Value = '  -1.67%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.93'
$ws.Range("E24").Value = '  -0.33%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.314'
$ws.Range("E25").Value = '  -1.49%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.86'
$ws.Range("E26").Value = '  -0.34%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.52'
$ws.Range("E27").Value = '  -2.57%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.127'
$ws.Range("E28").Value = '  +2.54%  '

# Row 29
$ws.Range("E29").Value = '  +0.13%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.13'
$ws.Range("E30").Value = '  -0.41%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.9736'
$ws.Range("E31").Value = '  +0.74%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09405'
$ws.Range("E32").Value = '  -1.00%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.618'
$ws.Range("E33").Value = '  -1.48%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.308'
$ws.Range("E34").Value = '  -0.66%  '

# Row 35
$ws.Range("E35").Value = '  -2.71%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02229'
$ws.Range("E36").Value = '  -0.72%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06010'
$ws.Range("E37").Value = '  -1.75%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.361'
$ws.Range("E38").Value = '  +3.43%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.188'
$ws.Range("E39").Value = '  -2.32%  '

# Row 40
$ws.Range("E40").Value = '  -0.74%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5919'
$ws.Range("E41").Value = '  -0.92%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.40'
$ws.Range("E42").Value = '  +0.80%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1866'
$ws.Range("E43").Value = '  -1.17%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.243'
$ws.Range("E44").Value = '  -2.27%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5581'
$ws.Range("E45").Value = '  -1.79%  '

# Row 46
$ws.Range("E46").Value = '  -0.06%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.909'
$ws.Range("E47").Value = '  -1.26%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06711'
$ws.Range("E48").Value = '  -3.12%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '111.14'
$ws.Range("E49").Value = '  -2.14%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.051'
$ws.Range("E50").Value = '  -1.68%  '

# Row 51
$ws.Range("E51").Value = '  -0.90%  '
